$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 8.494800000000005
$ws.Range("B13").Value = 5.613200000000003
$ws.Range("B16").Value = 9.009000000000007
$ws.Range("B18").Value = 5.357700000000001
$ws.Range("B20").Value = 5.6624
